# Generate Report for Handback
# Update the timestamp / status cells that the report-generation step refreshes.

$wb = $excel.ActiveWorkbook

# --- Overview sheet -------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("G3").Value = "2016-09-08 04:23:21"
$wsOverview.Range("G5").Value = "2016-09-08 04:23:21"

# --- zh-cn sheet ------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("E3").Value = "mt"
$wsZhCn.Range("E5").Value = "mt"
$wsZhCn.Range("H3").Value = "2016-09-08 04:23:15"
$wsZhCn.Range("H5").Value = "2016-09-08 04:23:15"
$wsZhCn.Range("K3").Value = "2016-09-08 04:23:35"
$wsZhCn.Range("K5").Value = "2016-09-08 04:23:35"

# --- de-de sheet ------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("E3").Value = "mt"
$wsDeDe.Range("E5").Value = "mt"
$wsDeDe.Range("H3").Value = "2016-09-08 04:23:21"
$wsDeDe.Range("H5").Value = "2016-09-08 04:23:21"
$wsDeDe.Range("K3").Value = "2016-09-08 04:23:43"
$wsDeDe.Range("K5").Value = "2016-09-08 04:23:43"
